$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add available meter types example: set D6:D10 to new meter type values
# that introduce three new shared strings: ЭЛ, ТЕПЛО, ГАЗ
$ws.Range("D6").Value = "ЭЛ"
$ws.Range("D7").Value = "ЭЛ"
$ws.Range("D8").Value = "ТЕПЛО"
$ws.Range("D9").Value = "ТЕПЛО"
$ws.Range("D10").Value = "ГАЗ"

# Move the selection cursor to D16 (matches saved cursor position)
$ws.Range("D16").Select()
